# Insert a new weekly data row at row 102 (pushing the existing rows 102..211
# down to 103..212), then populate the new row with the reported price data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(102).Insert()

$ws.Range("A102").Value = 9
$ws.Range("B102").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C102").Value = 'Metropolitana'
$ws.Range("D102").Value = 44539
$ws.Range("E102").Value = 13
$ws.Range("F102").Value = 100112001
$ws.Range("G102").Value = 'Berenjena'
$ws.Range("H102").Value = 'Sin especificar'
$ws.Range("I102").Value = 'Primera'
$ws.Range("J102").Value = 106
$ws.Range("K102").Value = 8000
$ws.Range("L102").Value = 9000
$ws.Range("M102").Value = 8500
$ws.Range("N102").Value = '$/caja 50 unidades'
$ws.Range("O102").Value = 'Región de Arica y Parinacota'
$ws.Range("P102").Value = 170
$ws.Range("Q102").Value = 50
$ws.Range("R102").Value = 'Hortaliza'
